$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.555934
$ws.Range("H2").Value = 1.667802
$ws.Range("I2").Value = 0.005745252779589096
$ws.Range("J2").Value = 0.005745252779589094
$ws.Range("M2").Value = 4.467713
$ws.Range("N2").Value = 13.403139
$ws.Range("O2").Value = 0.4045017685715556
$ws.Range("P2").Value = 0.4045017685715556
$ws.Range("Q2").Value = 2.483753558942
$ws.Range("R2").Value = 22.353782030478
$ws.Range("S2").Value = 0.002323964910234435
$ws.Range("T2").Value = 0.002323964910234434

$ws.Range("G3").Value = 0.555934
$ws.Range("H3").Value = 1.667802
$ws.Range("I3").Value = 0.005745252779589096
$ws.Range("J3").Value = 0.005745252779589094
$ws.Range("M3").Value = 3.429517666666667
$ws.Range("N3").Value = 10.288553
$ws.Range("O3").Value = 0.3105047171817127
$ws.Range("P3").Value = 0.3105047171817127
$ws.Range("Q3").Value = 1.906585474500667
$ws.Range("R3").Value = 17.159269270506
$ws.Range("S3").Value = 0.001783928089463761
$ws.Range("T3").Value = 0.001783928089463761

$ws.Range("G4").Value = 0.555934
$ws.Range("H4").Value = 1.667802
$ws.Range("I4").Value = 0.005745252779589096
$ws.Range("J4").Value = 0.005745252779589094
$ws.Range("M4").Value = 2.691692
$ws.Range("N4").Value = 8.075076000000001
$ws.Range("O4").Value = 0.2437028015116252
$ws.Range("P4").Value = 0.2437028015116253
$ws.Range("Q4").Value = 1.496403100328
$ws.Range("R4").Value = 13.467627902952
$ws.Range("S4").Value = 0.001400134197778314
$ws.Range("T4").Value = 0.001400134197778314

$ws.Range("G5").Value = 0.555934
$ws.Range("H5").Value = 1.667802
$ws.Range("I5").Value = 0.005745252779589096
$ws.Range("J5").Value = 0.005745252779589094
$ws.Range("M5").Value = 0.456055
$ws.Range("N5").Value = 1.368165
$ws.Range("O5").Value = 0.04129071273510649
$ws.Range("P5").Value = 0.04129071273510649
$ws.Range("Q5").Value = 0.25353648037
$ws.Range("R5").Value = 2.28182832333
$ws.Range("S5").Value = 0.0002372255821125854
$ws.Range("T5").Value = 0.0002372255821125854

$ws.Range("G6").Value = 79.68771233333334
$ws.Range("H6").Value = 239.063137
$ws.Range("I6").Value = 0.823525905561055
$ws.Range("J6").Value = 0.823525905561055
$ws.Range("M6").Value = 4.467713
$ws.Range("N6").Value = 13.403139
$ws.Range("O6").Value = 0.4045017685715556
$ws.Range("P6").Value = 0.4045017685715556
$ws.Range("Q6").Value = 356.0218283318937
$ws.Range("R6").Value = 3204.196454987043
$ws.Range("S6").Value = 0.3331176852639386
$ws.Range("T6").Value = 0.3331176852639386

$ws.Range("G7").Value = 79.68771233333334
$ws.Range("H7").Value = 239.063137
$ws.Range("I7").Value = 0.823525905561055
$ws.Range("J7").Value = 0.823525905561055
$ws.Range("M7").Value = 3.429517666666667
$ws.Range("N7").Value = 10.288553
$ws.Range("O7").Value = 0.3105047171817127
$ws.Range("P7").Value = 0.3105047171817127
$ws.Range("Q7").Value = 273.2904172634179
$ws.Range("R7").Value = 2459.613755370761
$ws.Range("S7").Value = 0.2557086783980492
$ws.Range("T7").Value = 0.2557086783980492

$ws.Range("G8").Value = 79.68771233333334
$ws.Range("H8").Value = 239.063137
$ws.Range("I8").Value = 0.823525905561055
$ws.Range("J8").Value = 0.823525905561055
$ws.Range("M8").Value = 2.691692
$ws.Range("N8").Value = 8.075076000000001
$ws.Range("O8").Value = 0.2437028015116252
$ws.Range("P8").Value = 0.2437028015116253
$ws.Range("Q8").Value = 214.4947777859347
$ws.Range("R8").Value = 1930.453000073412
$ws.Range("S8").Value = 0.2006955703026272
$ws.Range("T8").Value = 0.2006955703026272

$ws.Range("G9").Value = 79.68771233333334
$ws.Range("H9").Value = 239.063137
$ws.Range("I9").Value = 0.823525905561055
$ws.Range("J9").Value = 0.823525905561055
$ws.Range("M9").Value = 0.456055
$ws.Range("N9").Value = 1.368165
$ws.Range("O9").Value = 0.04129071273510649
$ws.Range("P9").Value = 0.04129071273510649
$ws.Range("Q9").Value = 36.34197964817834
$ws.Range("R9").Value = 327.077816833605
$ws.Range("S9").Value = 0.03400397159643996
$ws.Range("T9").Value = 0.03400397159643996

$ws.Range("G10").Value = 0.3446996666666666
$ws.Range("H10").Value = 1.034099
$ws.Range("I10").Value = 0.003562269474506148
$ws.Range("J10").Value = 0.003562269474506148
$ws.Range("M10").Value = 4.467713
$ws.Range("N10").Value = 13.403139
$ws.Range("O10").Value = 0.4045017685715556
$ws.Range("P10").Value = 0.4045017685715556
$ws.Range("Q10").Value = 1.540019181862333
$ws.Range("R10").Value = 13.860172636761
$ws.Range("S10").Value = 0.001440944302566203
$ws.Range("T10").Value = 0.001440944302566203

$ws.Range("G11").Value = 0.3446996666666666
$ws.Range("H11").Value = 1.034099
$ws.Range("I11").Value = 0.003562269474506148
$ws.Range("J11").Value = 0.003562269474506148
$ws.Range("M11").Value = 3.429517666666667
$ws.Range("N11").Value = 10.288553
$ws.Range("O11").Value = 0.3105047171817127
$ws.Range("P11").Value = 0.3105047171817127
$ws.Range("Q11").Value = 1.182153596527444
$ws.Range("R11").Value = 10.639382368747
$ws.Range("S11").Value = 0.00110610147570658
$ws.Range("T11").Value = 0.00110610147570658

$ws.Range("G12").Value = 0.3446996666666666
$ws.Range("H12").Value = 1.034099
$ws.Range("I12").Value = 0.003562269474506148
$ws.Range("J12").Value = 0.003562269474506148
$ws.Range("M12").Value = 2.691692
$ws.Range("N12").Value = 8.075076000000001
$ws.Range("O12").Value = 0.2437028015116252
$ws.Range("P12").Value = 0.2437028015116253
$ws.Range("Q12").Value = 0.9278253351693333
$ws.Range("R12").Value = 8.350428016524
$ws.Range("S12").Value = 0.0008681350506764934
$ws.Range("T12").Value = 0.0008681350506764934

$ws.Range("G13").Value = 0.3446996666666666
$ws.Range("H13").Value = 1.034099
$ws.Range("I13").Value = 0.003562269474506148
$ws.Range("J13").Value = 0.003562269474506148
$ws.Range("M13").Value = 0.456055
$ws.Range("N13").Value = 1.368165
$ws.Range("O13").Value = 0.04129071273510649
$ws.Range("P13").Value = 0.04129071273510649
$ws.Range("Q13").Value = 0.1572020064816667
$ws.Range("R13").Value = 1.414818058335
$ws.Range("S13").Value = 0.0001470886455568721
$ws.Range("T13").Value = 0.0001470886455568721

$ws.Range("G14").Value = 16.17571666666667
$ws.Range("H14").Value = 48.52715
$ws.Range("I14").Value = 0.1671665721848498
$ws.Range("J14").Value = 0.1671665721848498
$ws.Range("M14").Value = 4.467713
$ws.Range("N14").Value = 13.403139
$ws.Range("O14").Value = 0.4045017685715556
$ws.Range("P14").Value = 0.4045017685715556
$ws.Range("Q14").Value = 72.26845963598333
$ws.Range("R14").Value = 650.41613672385
$ws.Range("S14").Value = 0.06761917409481637
$ws.Range("T14").Value = 0.06761917409481637

$ws.Range("G15").Value = 16.17571666666667
$ws.Range("H15").Value = 48.52715
$ws.Range("I15").Value = 0.1671665721848498
$ws.Range("J15").Value = 0.1671665721848498
$ws.Range("M15").Value = 3.429517666666667
$ws.Range("N15").Value = 10.288553
$ws.Range("O15").Value = 0.3105047171817127
$ws.Range("P15").Value = 0.3105047171817127
$ws.Range("Q15").Value = 55.47490607932778
$ws.Range("R15").Value = 499.27415471395
$ws.Range("S15").Value = 0.05190600921849317
$ws.Range("T15").Value = 0.05190600921849316

$ws.Range("G16").Value = 16.17571666666667
$ws.Range("H16").Value = 48.52715
$ws.Range("I16").Value = 0.1671665721848498
$ws.Range("J16").Value = 0.1671665721848498
$ws.Range("M16").Value = 2.691692
$ws.Range("N16").Value = 8.075076000000001
$ws.Range("O16").Value = 0.2437028015116252
$ws.Range("P16").Value = 0.2437028015116253
$ws.Range("Q16").Value = 43.54004714593334
$ws.Range("R16").Value = 391.8604243134
$ws.Range("S16").Value = 0.04073896196054323
$ws.Range("T16").Value = 0.04073896196054323

$ws.Range("G17").Value = 16.17571666666667
$ws.Range("H17").Value = 48.52715
$ws.Range("I17").Value = 0.1671665721848498
$ws.Range("J17").Value = 0.1671665721848498
$ws.Range("M17").Value = 0.456055
$ws.Range("N17").Value = 1.368165
$ws.Range("O17").Value = 0.04129071273510649
$ws.Range("P17").Value = 0.04129071273510649
$ws.Range("Q17").Value = 7.377016464416667
$ws.Range("R17").Value = 66.39314817975
$ws.Range("S17").Value = 0.006902426910997078
$ws.Range("T17").Value = 0.006902426910997077

